$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.148.21'
$ws.Range("E2").Value = '  +3.64%  '

$ws.Range("D3").Value = '2.428.65'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("D5").Value = '''316.98'
$ws.Range("E5").Value = '  +3.37%  '

$ws.Range("D6").Value = '''102.71'
$ws.Range("E6").Value = '  +5.71%  '

$ws.Range("E7").Value = '  +1.75%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = '''0.526'
$ws.Range("E9").Value = '  +7.16%  '

$ws.Range("D10").Value = '''35.48'
$ws.Range("E10").Value = '  +1.50%  '

$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("E12").Value = '  -2.27%  '

$ws.Range("D13").Value = '''18.16'
$ws.Range("E13").Value = '  -1.76%  '

$ws.Range("D14").Value = '''7.03'
$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("D15").Value = '2.808.81'

$ws.Range("D16").Value = '2.452.08'
$ws.Range("E16").Value = '  +1.99%  '

$ws.Range("E17").Value = '  +2.11%  '

$ws.Range("D18").Value = '45.084.12'
$ws.Range("E18").Value = '  +3.42%  '

$ws.Range("D19").Value = '''12.28'
$ws.Range("E19").Value = '  +1.26%  '

$ws.Range("E20").Value = '  -0.95%  '

$ws.Range("E21").Value = '  +2.20%  '

$ws.Range("D22").Value = '''68.84'
$ws.Range("E22").Value = '  +0.63%  '

$ws.Range("D23").Value = '''244.17'
$ws.Range("E23").Value = '  +2.57%  '

$ws.Range("D24").Value = '''2.25'
$ws.Range("E24").Value = '  +0.56%  '

$ws.Range("D25").Value = '''2.50'
$ws.Range("E25").Value = '  +1.64%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").Value = '''25.48'
$ws.Range("E27").Value = '  +2.09%  '

$ws.Range("E28").Value = '  +1.47%  '

$ws.Range("E29").Value = '  -11.83%  '

$ws.Range("D30").Value = '''49.14'
$ws.Range("E30").Value = '  +2.70%  '

$ws.Range("D31").Value = '''32.95'
$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("D32").Value = '''20.32'
$ws.Range("E32").Value = '  +10.21%  '

$ws.Range("E33").Value = '  +5.57%  '

$ws.Range("D34").Value = '''5.21'
$ws.Range("E34").Value = '  +1.75%  '

$ws.Range("E35").Value = '  +0.15%  '

$ws.Range("D36").Value = '''0.0765'
$ws.Range("E36").Value = '  +1.86%  '

$ws.Range("E37").Value = '  -0.62%  '

$ws.Range("D38").Value = '''4.43'
$ws.Range("E38").Value = '  +0.78%  '

$ws.Range("E39").Value = '  -2.29%  '

$ws.Range("D40").Value = '''124.22'
$ws.Range("E40").Value = '  -6.73%  '

$ws.Range("E41").Value = '  -2.97%  '

$ws.Range("E42").Value = '  +0.96%  '

$ws.Range("D43").Value = '''20.57'
$ws.Range("E43").Value = '  -3.51%  '

$ws.Range("D44").Value = '''0.0288'
$ws.Range("E44").Value = '  +1.83%  '

$ws.Range("D45").Value = '1.933.42'
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("E46").Value = '  -3.00%  '

$ws.Range("E47").Value = '  +3.06%  '

$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '''1.81'
$ws.Range("E48").Value = '  +15.96%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''9.21'
$ws.Range("E49").Value = '  -0.95%  '

$ws.Range("D50").Value = '''76.66'
$ws.Range("E50").Value = '  +6.08%  '

$ws.Range("D51").Value = '''53.96'
$ws.Range("E51").Value = '  +2.49%  '
